$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy H1's formatting (bold, centered, bordered) to I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..15 : I = 1 (constant), J = same value as H
for ($r = 2; $r -le 15; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
